$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2450.2856
$ws.Range("I18").Value = 2450.2856
$ws.Range("K18").Value = 2450.2856
$ws.Range("M18").Value = -2166.2856
$ws.Range("H94").Value = 45575064
$ws.Range("I94").Value = 71437960
$ws.Range("K94").Value = 71437960
$ws.Range("M94").Value = -71437509
$ws.Range("H100").Value = 2952708.5
$ws.Range("I100").Value = 5106.8
$ws.Range("K100").Value = 5106.8
$ws.Range("M100").Value = -4565.8
$ws.Range("H137").Value = 3870604
$ws.Range("I137").Value = 860141.2
$ws.Range("J137").Value = 5441280
$ws.Range("K137").Value = 2580423.6
$ws.Range("L137").Value = 16323840
$ws.Range("M137").Value = -2577873.6
$ws.Range("N137").Value = -16328940

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 70000
$ws.Range("J24").Value = 70000
$ws.Range("L24").Value = 70000
$ws.Range("N24").Value = -70748
$ws.Range("H32").Value = 1580
$ws.Range("I32").Value = 1510.9678
$ws.Range("K32").Value = 1510.9678
$ws.Range("M32").Value = -1223.9678
$ws.Range("H61").Value = 5521.347
$ws.Range("I61").Value = 5671.15
$ws.Range("K61").Value = 5671.15
$ws.Range("M61").Value = -5459.15
$ws.Range("H74").Value = 12536917
$ws.Range("I74").Value = 16176109
$ws.Range("K74").Value = 16176109
$ws.Range("M74").Value = -16175235
$ws.Range("H77").Value = 12536917
$ws.Range("I77").Value = 16176109
$ws.Range("K77").Value = 80880545
$ws.Range("M77").Value = -80876177
$ws.Range("H100").Value = 70000
$ws.Range("J100").Value = 70000
$ws.Range("L100").Value = 70000
$ws.Range("N100").Value = -72164
$ws.Range("H102").Value = 5825.314
$ws.Range("I102").Value = 4078.0334
$ws.Range("K102").Value = 4078.0334
$ws.Range("M102").Value = -2456.0334
$ws.Range("H110").Value = 2631.4666
$ws.Range("J110").Value = 3356.4285
$ws.Range("L110").Value = 3356.4285
$ws.Range("N110").Value = -7446.4285
$ws.Range("H136").Value = 5521.347
$ws.Range("I136").Value = 5671.15
$ws.Range("K136").Value = 17013.45
$ws.Range("M136").Value = -14463.45

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 779.6667
$ws.Range("J22").Value = 779
$ws.Range("L22").Value = 779
$ws.Range("N22").Value = -1479
$ws.Range("H31").Value = 1947.8354
$ws.Range("I31").Value = 857.53845
$ws.Range("K31").Value = 857.53845
$ws.Range("M31").Value = -562.53845
$ws.Range("H34").Value = 1947.8354
$ws.Range("I34").Value = 857.53845
$ws.Range("K34").Value = 857.53845
$ws.Range("M34").Value = -655.53845
$ws.Range("H99").Value = 6253351
$ws.Range("I99").Value = 13891575
$ws.Range("J99").Value = 3894.7273
$ws.Range("K99").Value = 13891575
$ws.Range("L99").Value = 3894.7273
$ws.Range("M99").Value = -13890077
$ws.Range("N99").Value = -6890.7273
$ws.Range("H107").Value = 28608.084
$ws.Range("I107").Value = 47399.715
$ws.Range("K107").Value = 47399.715
$ws.Range("M107").Value = -45479.715
$ws.Range("H126").Value = 6253351
$ws.Range("I126").Value = 13891575
$ws.Range("J126").Value = 3894.7273
$ws.Range("K126").Value = 41674725
$ws.Range("L126").Value = 11684.1819
$ws.Range("M126").Value = -41672255
$ws.Range("N126").Value = -16624.1819
$ws.Range("H132").Value = 13355499
$ws.Range("I132").Value = 17551280
$ws.Range("K132").Value = 52653840
$ws.Range("M132").Value = -52651310
$ws.Range("H134").Value = 1694577.6
$ws.Range("I134").Value = 2506362.5
$ws.Range("K134").Value = 7519087.5
$ws.Range("M134").Value = -7516552.5
$ws.Range("H135").Value = 99998.664
$ws.Range("J135").Value = 99998.664
$ws.Range("L135").Value = 99998.664
$ws.Range("N135").Value = -110138.664

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 154.36363
$ws.Range("I29").Value = 71
$ws.Range("J29").Value = 376.66666
$ws.Range("K29").Value = 213
$ws.Range("L29").Value = 1129.99998
$ws.Range("M29").Value = 64
$ws.Range("N29").Value = -1683.99998
$ws.Range("H92").Value = 656.1579
$ws.Range("J92").Value = 842.44446
$ws.Range("L92").Value = 2527.33338
$ws.Range("N92").Value = -5023.33338
$ws.Range("H113").Value = 1289.8276
$ws.Range("J113").Value = 1579.421
$ws.Range("L113").Value = 4738.263
$ws.Range("N113").Value = -9078.262999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 29666.334
$ws.Range("J26").Value = 29666.334
$ws.Range("L26").Value = 29666.334
$ws.Range("N26").Value = -30226.334
$ws.Range("H50").Value = 29666.334
$ws.Range("J50").Value = 29666.334
$ws.Range("L50").Value = 29666.334
$ws.Range("N50").Value = -30662.334
$ws.Range("H97").Value = 6380.41
$ws.Range("I97").Value = 8686.962
$ws.Range("K97").Value = 8686.962
$ws.Range("M97").Value = -8190.962
$ws.Range("H102").Value = 12786.286
$ws.Range("I102").Value = 13539.154
$ws.Range("J102").Value = 2999
$ws.Range("K102").Value = 13539.154
$ws.Range("L102").Value = 2999
$ws.Range("M102").Value = -11917.154
$ws.Range("N102").Value = -6243
$ws.Range("H122").Value = 7739.64
$ws.Range("I122").Value = 8950.200000000001
$ws.Range("K122").Value = 26850.6
$ws.Range("M122").Value = -24400.6
$ws.Range("H132").Value = 4546.9375
$ws.Range("I132").Value = 3665.7932
$ws.Range("K132").Value = 10997.3796
$ws.Range("M132").Value = -8467.3796
$ws.Range("H135").Value = 115995
$ws.Range("J135").Value = 115995
$ws.Range("L135").Value = 115995
$ws.Range("N135").Value = -126135

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8052.484
$ws.Range("I22").Value = 10879.8125
$ws.Range("J22").Value = 5036.6665
$ws.Range("K22").Value = 10879.8125
$ws.Range("L22").Value = 5036.6665
$ws.Range("M22").Value = -10584.8125
$ws.Range("N22").Value = -5626.6665
$ws.Range("H27").Value = 8052.484
$ws.Range("I27").Value = 10879.8125
$ws.Range("J27").Value = 5036.6665
$ws.Range("K27").Value = 10879.8125
$ws.Range("L27").Value = 5036.6665
$ws.Range("M27").Value = -10772.8125
$ws.Range("N27").Value = -5250.6665
$ws.Range("H40").Value = 46738.223
$ws.Range("I40").Value = 60808.25
$ws.Range("K40").Value = 60808.25
$ws.Range("M40").Value = -60672.25
$ws.Range("H61").Value = 3095.25
$ws.Range("I61").Value = 2524.3333
$ws.Range("J61").Value = 3666.1667
$ws.Range("K61").Value = 2524.3333
$ws.Range("L61").Value = 3666.1667
$ws.Range("M61").Value = -2322.3333
$ws.Range("N61").Value = -4070.1667
$ws.Range("H113").Value = 3095.25
$ws.Range("I113").Value = 2524.3333
$ws.Range("J113").Value = 3666.1667
$ws.Range("K113").Value = 2524.3333
$ws.Range("L113").Value = 3666.1667
$ws.Range("M113").Value = -354.3332999999998
$ws.Range("N113").Value = -8006.1667
$ws.Range("H127").Value = 167166670
$ws.Range("J127").Value = 599999.2
$ws.Range("L127").Value = 599999.2
$ws.Range("N127").Value = -609919.2
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H136").Value = 4701.0244
$ws.Range("I136").Value = 1569
$ws.Range("J136").Value = 8703.056
$ws.Range("K136").Value = 4707
$ws.Range("L136").Value = 26109.168
$ws.Range("M136").Value = -2157
$ws.Range("N136").Value = -31209.168

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16864.334
$ws.Range("I132").Value = 19970
$ws.Range("K132").Value = 59910
$ws.Range("M132").Value = -57380
$ws.Range("H136").Value = 4098.3145
$ws.Range("I136").Value = 3378.56
$ws.Range("J136").Value = 5897.7
$ws.Range("K136").Value = 10135.68
$ws.Range("L136").Value = 17693.1
$ws.Range("M136").Value = -7585.68
$ws.Range("N136").Value = -22793.1

Write-Host "edits applied"